$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "2025-08-13 04:01:50 UTC"
$ws.Range("B6").Value = "2025-08-13 09:31:50 IST"
$ws.Range("C6").Value = "SKIPPED"
$ws.Range("D6").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E6").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = ""

$ws.Range("A6:H6").HorizontalAlignment = -4108
$ws.Range("A6:H6").VerticalAlignment = -4108
